$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '38.780.72'
$ws.Range("E2").Value = '  +1.20%  '
$ws.Range("D3").Value = '2.096.35'
$ws.Range("E3").Value = '  +0.91%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '226.83'
$ws.Range("E5").Value = '  -0.55%  '
$ws.Range("D6").Value = '0.616'
$ws.Range("E6").Value = '  +1.13%  '
$ws.Range("D7").Value = '61.87'
$ws.Range("E7").Value = '  +3.00%  '
$ws.Range("E8").Value = '  -0.07%  '
$ws.Range("D9").Value = '0.387'
$ws.Range("E9").Value = '  +1.82%  '
$ws.Range("D10").Value = '0.0839'
$ws.Range("E10").Value = '  +0.92%  '
$ws.Range("E11").Value = '  -0.23%  '
$ws.Range("D12").Value = '15.71'
$ws.Range("E12").Value = '  +6.31%  '
$ws.Range("D13").Value = '2.413.95'
$ws.Range("E13").Value = '  +1.02%  '
$ws.Range("D14").Value = '21.88'
$ws.Range("E14").Value = '  -1.59%  '
$ws.Range("D15").Value = '0.802'
$ws.Range("E15").Value = '  +3.06%  '
$ws.Range("D16").Value = '5.49'
$ws.Range("E16").Value = '  +1.11%  '
$ws.Range("D17").Value = '2.151.01'
$ws.Range("E17").Value = '  +3.53%  '
$ws.Range("D18").Value = '38.736.12'
$ws.Range("E18").Value = '  +1.10%  '
$ws.Range("D19").Value = '71.56'
$ws.Range("E19").Value = '  +0.15%  '
$ws.Range("D20").Value = '6.06'
$ws.Range("E20").Value = '  +1.13%  '
$ws.Range("D21").Value = '0.0₃0844'
$ws.Range("E21").Value = '  +1.77%  '
$ws.Range("D22").Value = '227.35'
$ws.Range("E22").Value = '  +1.27%  '
$ws.Range("E23").Value = '  -0.04%  '
$ws.Range("E24").Value = '  -2.68%  '
$ws.Range("D25").Value = '2.31'
$ws.Range("E25").Value = '  -0.42%  '
$ws.Range("D26").Value = '9.65'
$ws.Range("E26").Value = '  +3.16%  '
$ws.Range("D27").Value = '170.99'
$ws.Range("E27").Value = '  +0.79%  '
$ws.Range("D28").Value = '0.135'
$ws.Range("E28").Value = '  -0.74%  '
$ws.Range("D29").Value = '1.41'
$ws.Range("E29").Value = '  +3.92%  '
$ws.Range("D30").Value = '19.29'
$ws.Range("E30").Value = '  +1.64%  '
$ws.Range("D31").Value = '2.54'
$ws.Range("E31").Value = '  +9.77%  '
$ws.Range("D32").Value = '0.120'
$ws.Range("E32").Value = '  +0.56%  '
$ws.Range("B33").Value = 'THORChain'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D33").Value = '7.21'
$ws.Range("E33").Value = '  +13.85%  '
$ws.Range("B34").Value = 'Filecoin'
$ws.Range("C34").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D34").Value = '4.56'
$ws.Range("E34").Value = '  +1.70%  '
$ws.Range("D35").Value = '4.72'
$ws.Range("E35").Value = '  -1.35%  '
$ws.Range("D36").Value = '0.0614'
$ws.Range("E36").Value = '  +1.70%  '
$ws.Range("E37").Value = '  +0.48%  '
$ws.Range("D38").Value = '3.51'
$ws.Range("E38").Value = '  +0.12%  '
$ws.Range("E39").Value = '  +0.03%  '
$ws.Range("D40").Value = '17.93'
$ws.Range("E40").Value = '  -0.94%  '
$ws.Range("E41").Value = '  +3.59%  '
$ws.Range("D42").Value = '101.38'
$ws.Range("E42").Value = '  +1.29%  '
$ws.Range("D43").Value = '1.523.45'
$ws.Range("E43").Value = '  -0.94%  '
$ws.Range("D44").Value = '1.19'
$ws.Range("E44").Value = '  +7.40%  '
$ws.Range("E45").Value = '  +0.31%  '
$ws.Range("D46").Value = '7.78'
$ws.Range("E46").Value = '  +2.17%  '
$ws.Range("D47").Value = '0.0910'
$ws.Range("E47").Value = '  -0.98%  '
$ws.Range("D48").Value = '1.08'
$ws.Range("E48").Value = '  +5.43%  '
$ws.Range("D49").Value = '4.14'
$ws.Range("E49").Value = '  +2.01%  '
$ws.Range("D50").Value = '2.96'
$ws.Range("E50").Value = '  -0.19%  '
$ws.Range("D51").Value = '2.300.70'
$ws.Range("E51").Value = '  +0.99%  '
